$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Lithuania A Lyga results base was updated: rows were re-sorted/re-matched
# (fixture IDs swapped between adjacent rows), so we rewrite the data-bearing
# columns (B, F-AC) for each affected row to reflect the corrected pairing.
# Columns A (running index), C (Div), D (Div Original Name) and E (Date) are
# unchanged for every affected row.

# Row 26
$ws.Range("B26").Value = 6732773
$ws.Range("F26").Value = 'Suduva Marijampole'
$ws.Range("G26").Value = 'Hegelmann Litauen'
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 'A'
$ws.Range("K26").Value = 5
$ws.Range("L26").Value = 3.8
$ws.Range("M26").Value = 1.533
$ws.Range("N26").Value = 5
$ws.Range("O26").Value = 4.2
$ws.Range("P26").Value = 1.533
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = 1.875
$ws.Range("S26").Value = 1.925
$ws.Range("T26").Value = 2.5
$ws.Range("U26").Value = 1.9
$ws.Range("V26").Value = 1.9
$ws.Range("W26").Value = -1
$ws.Range("X26").Value = -1
$ws.Range("Y26").Value = 0.5329999999999999
$ws.Range("Z26").Value = 0
$ws.Range("AA26").Value = -0
$ws.Range("AB26").Value = -1
$ws.Range("AC26").Value = 0.8999999999999999

# Row 27
$ws.Range("B27").Value = 6732711
$ws.Range("F27").Value = 'Banga Gargzdai'
$ws.Range("G27").Value = 'FK Zalgiris Vilnius'
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 'A'
$ws.Range("K27").Value = 5
$ws.Range("L27").Value = 3.6
$ws.Range("M27").Value = 1.571
$ws.Range("N27").Value = 11
$ws.Range("O27").Value = 4.75
$ws.Range("P27").Value = 1.25
$ws.Range("Q27").Value = 1.5
$ws.Range("R27").Value = 1.975
$ws.Range("S27").Value = 1.825
$ws.Range("T27").Value = 2.5
$ws.Range("U27").Value = 1.8
$ws.Range("V27").Value = 2
$ws.Range("W27").Value = -1
$ws.Range("X27").Value = -1
$ws.Range("Y27").Value = 0.25
$ws.Range("Z27").Value = -1
$ws.Range("AA27").Value = 0.825
$ws.Range("AB27").Value = 0.8
$ws.Range("AC27").Value = -1

# Row 89
$ws.Range("B89").Value = 6732827
$ws.Range("F89").Value = 'FK Dziugas Telsiai'
$ws.Range("G89").Value = 'FK Kauno Zalgiris'
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 2
$ws.Range("J89").Value = 'A'
$ws.Range("K89").Value = 6
$ws.Range("L89").Value = 3.9
$ws.Range("M89").Value = 1.444
$ws.Range("N89").Value = 4.75
$ws.Range("O89").Value = 3.6
$ws.Range("P89").Value = 1.65
$ws.Range("Q89").Value = 0.75
$ws.Range("R89").Value = 1.9
$ws.Range("S89").Value = 1.9
$ws.Range("T89").Value = 2.5
$ws.Range("U89").Value = 1.95
$ws.Range("V89").Value = 1.85
$ws.Range("W89").Value = -1
$ws.Range("X89").Value = -1
$ws.Range("Y89").Value = 0.6499999999999999
$ws.Range("Z89").Value = -1
$ws.Range("AA89").Value = 0.8999999999999999
$ws.Range("AB89").Value = -1
$ws.Range("AC89").Value = 0.8500000000000001

# Row 90
$ws.Range("B90").Value = 7326568
$ws.Range("F90").Value = 'Hegelmann Litauen'
$ws.Range("G90").Value = 'Panevezys'
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 'D'
$ws.Range("K90").Value = 2.375
$ws.Range("L90").Value = 3.2
$ws.Range("M90").Value = 2.625
$ws.Range("N90").Value = 2.7
$ws.Range("O90").Value = 3.2
$ws.Range("P90").Value = 2.3
$ws.Range("Q90").Value = 0
$ws.Range("R90").Value = 2.05
$ws.Range("S90").Value = 1.75
$ws.Range("T90").Value = 2.25
$ws.Range("U90").Value = 1.875
$ws.Range("V90").Value = 1.925
$ws.Range("W90").Value = -1
$ws.Range("X90").Value = 2.2
$ws.Range("Y90").Value = -1
$ws.Range("Z90").Value = 0
$ws.Range("AA90").Value = -0
$ws.Range("AB90").Value = -1
$ws.Range("AC90").Value = 0.925

# Row 100
$ws.Range("B100").Value = 6732727
$ws.Range("F100").Value = 'FK Zalgiris Vilnius'
$ws.Range("G100").Value = 'FK Dainava Alytus'
$ws.Range("H100").Value = 1
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 'H'
$ws.Range("K100").Value = 1.285
$ws.Range("L100").Value = 5.5
$ws.Range("M100").Value = 6.5
$ws.Range("N100").Value = 1.3
$ws.Range("O100").Value = 5.5
$ws.Range("P100").Value = 6
$ws.Range("Q100").Value = -1.5
$ws.Range("R100").Value = 1.9
$ws.Range("S100").Value = 1.9
$ws.Range("T100").Value = 2.75
$ws.Range("U100").Value = 1.8
$ws.Range("V100").Value = 2
$ws.Range("W100").Value = 0.3
$ws.Range("X100").Value = -1
$ws.Range("Y100").Value = -1
$ws.Range("Z100").Value = -1
$ws.Range("AA100").Value = 0.8999999999999999
$ws.Range("AB100").Value = -1
$ws.Range("AC100").Value = 1

# Row 101
$ws.Range("B101").Value = 6732836
$ws.Range("F101").Value = 'FK Siauliai'
$ws.Range("G101").Value = 'Banga Gargzdai'
$ws.Range("H101").Value = 3
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 'H'
$ws.Range("K101").Value = 1.222
$ws.Range("L101").Value = 5.5
$ws.Range("M101").Value = 9
$ws.Range("N101").Value = 1.363
$ws.Range("O101").Value = 4.5
$ws.Range("P101").Value = 7
$ws.Range("Q101").Value = -1.25
$ws.Range("R101").Value = 1.9
$ws.Range("S101").Value = 1.9
$ws.Range("T101").Value = 2.5
$ws.Range("U101").Value = 1.975
$ws.Range("V101").Value = 1.825
$ws.Range("W101").Value = 0.363
$ws.Range("X101").Value = -1
$ws.Range("Y101").Value = -1
$ws.Range("Z101").Value = 0.8999999999999999
$ws.Range("AA101").Value = -1
$ws.Range("AB101").Value = 0.9750000000000001
$ws.Range("AC101").Value = -1

# Row 102
$ws.Range("B102").Value = 7465686
$ws.Range("F102").Value = 'FK Kauno Zalgiris'
$ws.Range("G102").Value = 'Hegelmann Litauen'
$ws.Range("H102").Value = 4
$ws.Range("I102").Value = 2
$ws.Range("J102").Value = 'H'
$ws.Range("K102").Value = 2.3
$ws.Range("L102").Value = 4
$ws.Range("M102").Value = 2.3
$ws.Range("N102").Value = 2.55
$ws.Range("O102").Value = 4
$ws.Range("P102").Value = 2.2
$ws.Range("Q102").Value = 0.25
$ws.Range("R102").Value = 1.8
$ws.Range("S102").Value = 2
$ws.Range("T102").Value = 2.75
$ws.Range("U102").Value = 1.85
$ws.Range("V102").Value = 1.95
$ws.Range("W102").Value = 1.55
$ws.Range("X102").Value = -1
$ws.Range("Y102").Value = -1
$ws.Range("Z102").Value = 0.8
$ws.Range("AA102").Value = -1
$ws.Range("AB102").Value = 0.8500000000000001
$ws.Range("AC102").Value = -1

# Row 104
$ws.Range("B104").Value = 6732834
$ws.Range("F104").Value = 'Panevezys'
$ws.Range("G104").Value = 'FK Dziugas Telsiai'
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 'D'
$ws.Range("K104").Value = 1.25
$ws.Range("L104").Value = 5.5
$ws.Range("M104").Value = 7.5
$ws.Range("N104").Value = 1.45
$ws.Range("O104").Value = 4.5
$ws.Range("P104").Value = 5
$ws.Range("Q104").Value = -1
$ws.Range("R104").Value = 1.775
$ws.Range("S104").Value = 2.025
$ws.Range("T104").Value = 2.5
$ws.Range("U104").Value = 1.875
$ws.Range("V104").Value = 1.925
$ws.Range("W104").Value = -1
$ws.Range("X104").Value = 3.5
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = -1
$ws.Range("AA104").Value = 1.025
$ws.Range("AB104").Value = -1
$ws.Range("AC104").Value = 0.925
